# Generate Report for Handback
# Fills in the handback information (status / target file / handback file /
# handback datetime) for the two rows ("ce4e7e07..." and "fa6de6f8...")
# that were still "Ready for handoff" before this run, on all three sheets:
# Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: columns E (zh-cn status) / F (de-de status) for the two
# rows that just got handed back.
# ---------------------------------------------------------------------
$ws1.Range("E4").Value = $handedBack
$ws1.Range("F4").Value = $handedBack
$ws1.Range("E5").Value = $handedBack
$ws1.Range("F5").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet: Status, Latest Target File (+hyperlink), Latest Handback
# File, Latest Handback DateTime.
# ---------------------------------------------------------------------
$ws2.Range("C4").Value = $handedBack
$ws2.Range("J4").Value = "ce4e7e07-7b4d-4fc3-959a-856f99577c2a.9ccf67ceaf2a166d132ffddb110059c3b601f481.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-08-17 08:29:33"
$ws2.Hyperlinks.Add($ws2.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d018da123e40e3e8df016c79127cc0e814f49b6/e2e/ce4e7e07-7b4d-4fc3-959a-856f99577c2a.md", $null, $null, "ce4e7e07-7b4d-4fc3-959a-856f99577c2a.md") | Out-Null

$ws2.Range("C5").Value = $handedBack
$ws2.Range("J5").Value = "fa6de6f8-b92d-4bd3-93f5-ced2971aa743.2ecfb03ad718204e1ccdb0289d17ae6344536df1.zh-cn.xlf"
$ws2.Range("K5").Value = "2016-08-17 08:29:33"
$ws2.Hyperlinks.Add($ws2.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d018da123e40e3e8df016c79127cc0e814f49b6/e2e/fa6de6f8-b92d-4bd3-93f5-ced2971aa743.md", $null, $null, "fa6de6f8-b92d-4bd3-93f5-ced2971aa743.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: same shape as zh-cn, different handoff xlf names / time.
# ---------------------------------------------------------------------
$ws3.Range("C4").Value = $handedBack
$ws3.Range("J4").Value = "ce4e7e07-7b4d-4fc3-959a-856f99577c2a.9ccf67ceaf2a166d132ffddb110059c3b601f481.de-de.xlf"
$ws3.Range("K4").Value = "2016-08-17 08:29:40"
$ws3.Hyperlinks.Add($ws3.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d018da123e40e3e8df016c79127cc0e814f49b6/e2e/ce4e7e07-7b4d-4fc3-959a-856f99577c2a.md", $null, $null, "ce4e7e07-7b4d-4fc3-959a-856f99577c2a.md") | Out-Null

$ws3.Range("C5").Value = $handedBack
$ws3.Range("J5").Value = "fa6de6f8-b92d-4bd3-93f5-ced2971aa743.2ecfb03ad718204e1ccdb0289d17ae6344536df1.de-de.xlf"
$ws3.Range("K5").Value = "2016-08-17 08:29:40"
$ws3.Hyperlinks.Add($ws3.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d018da123e40e3e8df016c79127cc0e814f49b6/e2e/fa6de6f8-b92d-4bd3-93f5-ced2971aa743.md", $null, $null, "fa6de6f8-b92d-4bd3-93f5-ced2971aa743.md") | Out-Null
